$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 905.3043
$ws.Range("I129").Value = 681.3333
$ws.Range("J129").Value = 920.93024
$ws.Range("K129").Value = 2043.9999
$ws.Range("L129").Value = 2762.79072
$ws.Range("M129").Value = 2956.0001
$ws.Range("N129").Value = -12762.79072

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1200
$ws.Range("I25").Value = 1200
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1200
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -798
$ws.Range("N25").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11005.811
$ws.Range("I32").Value = 7621.163
$ws.Range("J32").Value = 29433.334
$ws.Range("K32").Value = 7621.163
$ws.Range("L32").Value = 29433.334
$ws.Range("M32").Value = -7334.163
$ws.Range("N32").Value = -30007.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3004719.5
$ws.Range("I61").Value = 3473888.2
$ws.Range("J61").Value = 2040
$ws.Range("K61").Value = 3473888.2
$ws.Range("L61").Value = 2040
$ws.Range("M61").Value = -3473676.2
$ws.Range("N61").Value = -2464

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3004719.5
$ws.Range("I136").Value = 3473888.2
$ws.Range("J136").Value = 2040
$ws.Range("K136").Value = 10421664.6
$ws.Range("L136").Value = 6120
$ws.Range("M136").Value = -10419114.6
$ws.Range("N136").Value = -11220

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1578.8667
$ws.Range("I31").Value = 1389.0834
$ws.Range("J31").Value = 2338
$ws.Range("K31").Value = 1389.0834
$ws.Range("L31").Value = 2338
$ws.Range("M31").Value = -1094.0834
$ws.Range("N31").Value = -2928

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1578.8667
$ws.Range("I34").Value = 1389.0834
$ws.Range("J34").Value = 2338
$ws.Range("K34").Value = 1389.0834
$ws.Range("L34").Value = 2338
$ws.Range("M34").Value = -1187.0834
$ws.Range("N34").Value = -2742

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1388.5
$ws.Range("I122").Value = 777
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 2331
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = 119
$ws.Range("N122").Value = -10900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2292.1482
$ws.Range("I134").Value = 2146.6667
$ws.Range("J134").Value = 2801.3333
$ws.Range("K134").Value = 6440.000100000001
$ws.Range("L134").Value = 8403.999899999999
$ws.Range("M134").Value = -3905.000100000001
$ws.Range("N134").Value = -13473.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 41972
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 41972
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 41972
$ws.Range("N138").Value = -52252
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H139").Value = 48008.75
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 48008.75
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 48008.75
$ws.Range("N139").Value = -58288.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 49640
$ws.Range("I140").Value = 50000
$ws.Range("J140").Value = 49280
$ws.Range("K140").Value = 50000
$ws.Range("L140").Value = 49280
$ws.Range("M140").Value = -44820
$ws.Range("N140").Value = -59640

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 13080.4
$ws.Range("I16").Value = 15600.5
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 46801.5
$ws.Range("L16").Value = 9000
$ws.Range("M16").Value = -46628.5
$ws.Range("N16").Value = -9346

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 490
$ws.Range("I51").Value = 490
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 1470
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -1010

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1995.16
$ws.Range("I68").Value = 572.5
$ws.Range("J68").Value = 2664.647
$ws.Range("K68").Value = 1717.5
$ws.Range("L68").Value = 7993.941
$ws.Range("M68").Value = -906.5
$ws.Range("N68").Value = -9615.940999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1995.16
$ws.Range("I71").Value = 572.5
$ws.Range("J71").Value = 2664.647
$ws.Range("K71").Value = 5152.5
$ws.Range("L71").Value = 23981.823
$ws.Range("M71").Value = -1096.5
$ws.Range("N71").Value = -32093.823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 4600
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 4600
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 13800
$ws.Range("N75").Value = -15796
$ws.Range("M75").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 4600
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 4600
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 41400
$ws.Range("N78").Value = -51384
$ws.Range("M78").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 18107640
$ws.Range("I80").Value = 24143186
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 72429558
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -72428622
$ws.Range("N80").Value = -4872

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 18107640
$ws.Range("I83").Value = 24143186
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 217288674
$ws.Range("L83").Value = 9000
$ws.Range("M83").Value = -217283994
$ws.Range("N83").Value = -18360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1930.8173
$ws.Range("I131").Value = 4289.846
$ws.Range("J131").Value = 1547.475
$ws.Range("K131").Value = 12869.538
$ws.Range("L131").Value = 4642.424999999999
$ws.Range("M131").Value = -7829.537999999999
$ws.Range("N131").Value = -14722.425

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 22078.666
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 22078.666
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 22078.666
$ws.Range("N15").Value = -22654.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 22078.666
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 22078.666
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 22078.666
$ws.Range("N81").Value = -24074.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 22078.666
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 22078.666
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 66235.99800000001
$ws.Range("N84").Value = -76219.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 13523.527
$ws.Range("I122").Value = 2531.7222
$ws.Range("J122").Value = 24515.334
$ws.Range("K122").Value = 7595.1666
$ws.Range("L122").Value = 73546.00199999999
$ws.Range("M122").Value = -5145.1666
$ws.Range("N122").Value = -78446.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2235.75
$ws.Range("I132").Value = 1911.4286
$ws.Range("J132").Value = 2992.5
$ws.Range("K132").Value = 5734.2858
$ws.Range("L132").Value = 8977.5
$ws.Range("M132").Value = -3204.2858
$ws.Range("N132").Value = -14037.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("N48").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 13998.667
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 13998.667
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 13998.667
$ws.Range("N54").Value = -15286.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1401.9642
$ws.Range("I82").Value = 1198.2307
$ws.Range("J82").Value = 1578.5333
$ws.Range("K82").Value = 1198.2307
$ws.Range("L82").Value = 1578.5333
$ws.Range("M82").Value = -837.2307000000001
$ws.Range("N82").Value = -2300.5333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1401.9642
$ws.Range("I85").Value = 1198.2307
$ws.Range("J85").Value = 1578.5333
$ws.Range("K85").Value = 1198.2307
$ws.Range("L85").Value = 1578.5333
$ws.Range("M85").Value = 49.76929999999993
$ws.Range("N85").Value = -4074.5333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2258.8635
$ws.Range("I122").Value = 2042.0588
$ws.Range("J122").Value = 2996
$ws.Range("K122").Value = 6126.1764
$ws.Range("L122").Value = 8988
$ws.Range("M122").Value = -3676.1764
$ws.Range("N122").Value = -13888

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11024.23
$ws.Range("I132").Value = 13768.556
$ws.Range("J132").Value = 4849.5
$ws.Range("K132").Value = 41305.66800000001
$ws.Range("L132").Value = 14548.5
$ws.Range("M132").Value = -38775.66800000001
$ws.Range("N132").Value = -19608.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 26500
$ws.Range("I75").Value = 26000
$ws.Range("J75").Value = 27000
$ws.Range("K75").Value = 26000
$ws.Range("L75").Value = 27000
$ws.Range("M75").Value = -25064
$ws.Range("N75").Value = -28872

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H78").Value = 26500
$ws.Range("I78").Value = 26000
$ws.Range("J78").Value = 27000
$ws.Range("K78").Value = 78000
$ws.Range("L78").Value = 81000
$ws.Range("M78").Value = -73320
$ws.Range("N78").Value = -90360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1520.3667
$ws.Range("I122").Value = 1520.3667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4561.1001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2111.1001
$ws.Range("N122").ClearContents()
